$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - styled like the other header cells (bold/border/center)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Data cells - time_taken values for each row
$ws.Range("F2").Value = "2021-10-05 13:41:03.955053"
$ws.Range("F3").Value = "2021-10-05 13:41:03.955064"
$ws.Range("F4").Value = "2021-10-05 13:41:03.955068"
$ws.Range("F5").Value = "2021-10-05 13:41:03.955071"
$ws.Range("F6").Value = "2021-10-05 13:41:03.955074"
